$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their original text storage (avoid Excel
# auto-converting numeric-looking strings like "1.41" into real numbers).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "92.089.71"
$ws.Range("E2").Value = "  -1.88%  "

$ws.Range("D3").Value = "3.334.53"
$ws.Range("E3").Value = "  -2.53%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "231.19"
$ws.Range("E5").Value = "  -2.08%  "

$ws.Range("D6").Value = "616.01"
$ws.Range("E6").Value = "  -3.38%  "

$ws.Range("D7").Value = "1.41"
$ws.Range("E7").Value = "  -1.12%  "

$ws.Range("E8").Value = "  -1.75%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").Value = "0.962"
$ws.Range("E10").Value = "  +0.55%  "

$ws.Range("D11").Value = "3.333.69"
$ws.Range("E11").Value = "  -2.62%  "

$ws.Range("D12").Value = "43.14"
$ws.Range("E12").Value = "  +3.04%  "

$ws.Range("E13").Value = "  -0.63%  "

$ws.Range("D14").Value = "6.16"
$ws.Range("E14").Value = "  +1.08%  "

$ws.Range("D15").Value = "91.947.38"
$ws.Range("E15").Value = "  -1.56%  "

$ws.Range("D16").Value = "3.959.42"
$ws.Range("E16").Value = "  -2.27%  "

$ws.Range("E17").Value = "  -2.55%  "

$ws.Range("E18").Value = "  -2.26%  "

$ws.Range("D19").Value = "3.330.40"
$ws.Range("E19").Value = "  -2.04%  "

$ws.Range("D20").Value = "17.45"
$ws.Range("E20").Value = "  -0.48%  "

$ws.Range("D21").Value = "10.91"
$ws.Range("E21").Value = "  -3.00%  "

$ws.Range("D22").Value = "3.46"
$ws.Range("E22").Value = "  +8.36%  "

$ws.Range("D23").Value = "494.44"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").Value = "0.443"
$ws.Range("E24").Value = "  -9.19%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("E26").Value = "  -2.52%  "

$ws.Range("D27").Value = "89.97"
$ws.Range("E27").Value = "  -0.56%  "

$ws.Range("D28").Value = "11.92"
$ws.Range("E28").Value = "  -0.33%  "

$ws.Range("D29").Value = "3.514.97"
$ws.Range("E29").Value = "  -2.00%  "

$ws.Range("E30").Value = "  -0.68%  "

$ws.Range("D31").Value = "11.19"
$ws.Range("E31").Value = "  -3.92%  "

$ws.Range("E32").Value = "  +2.47%  "

$ws.Range("E33").Value = "  -3.47%  "

$ws.Range("D34").Value = "0.994"
$ws.Range("E34").Value = "  -0.59%  "

$ws.Range("D35").Value = "0.173"
$ws.Range("E35").Value = "  -3.30%  "

$ws.Range("D36").Value = "28.41"
$ws.Range("E36").Value = "  -5.12%  "

$ws.Range("E37").Value = "  -4.59%  "

$ws.Range("D38").Value = "564.60"
$ws.Range("E38").Value = "  +3.93%  "

$ws.Range("E39").Value = "  -2.07%  "

$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("E41").Value = "  -0.39%  "

$ws.Range("E42").Value = "  -4.50%  "

$ws.Range("D43").Value = "0.870"
$ws.Range("E43").Value = "  -6.01%  "

$ws.Range("E44").Value = "  -1.48%  "

$ws.Range("D45").Value = "1.68"
$ws.Range("E45").Value = "  +0.49%  "

$ws.Range("D48").Value = "5.42"
$ws.Range("E48").Value = "  -1.78%  "

$ws.Range("D49").Value = "2.12"
$ws.Range("E49").Value = "  -0.52%  "

$ws.Range("E50").Value = "  +0.21%  "

$ws.Range("D51").Value = "51.57"
$ws.Range("E51").Value = "  -2.18%  "

# Row 46 and 47 swapped coins (VeChain <-> MantraDAO) with updated data
$ws.Range("B46").Value = "MantraDAO"
$ws.Range("C46").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D46").Value = "3.61"
$ws.Range("E46").Value = "  +5.45%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0414"
$ws.Range("E47").Value = "  +1.23%  "
